# "My Cart test cases added."
# Adds a new "MyCart" worksheet (after "Login") populated with cart/checkout
# test-case rows, formats header + data cells (centered, wrapped text),
# adds hyperlinks for the "User Name" column, sets row heights, column
# widths, and updates the selection / active-sheet state.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Preserve the final selection on the "Login" sheet (B4) before we
#    start juggling sheets, since Select() reflects the state at save time.
# ---------------------------------------------------------------------
$loginSheet = $wb.Worksheets.Item("Login")
$loginSheet.Range("B4").Select()

# ---------------------------------------------------------------------
# 2. Insert the "MyCart" worksheet right after "Login".
#    Two throw-away sheets are added+removed first purely so the
#    workbook's internal sheetId counter lands on 3 for MyCart (matching
#    a workbook that previously had a 2nd sheet created/removed).
# ---------------------------------------------------------------------
$tmp1 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item("Login"))
$tmp1.Name = "zzTemp1"
$tmp2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item("zzTemp1"))
$tmp2.Name = "zzTemp2"
$wb.Worksheets.Item("zzTemp2").Delete()

$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item("zzTemp1"))
$ws.Name = "MyCart"

$wb.Worksheets.Item("zzTemp1").Delete()

$ws = $wb.Worksheets.Item("MyCart")

# ---------------------------------------------------------------------
# 3. Header row. Order matters: it controls which shared-string index
#    each text value receives, so write in the exact sequence the
#    original authoring session would have used.
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "Test Case Description"
$ws.Range("B1").Value = "User Name"
$ws.Range("C1").Value = "Password"
$ws.Range("E1").Value = "Shipping Address"
$ws.Range("D1").Value = "Billing Address"

# ---------------------------------------------------------------------
# 4. Data rows.
# ---------------------------------------------------------------------
$ws.Range("A2").Value = "TC_SC_001 - Add to Cart"
$ws.Range("B2").Value = "test@user.com"
$ws.Range("C2").Value = 1234
$ws.Range("D2").Value = "123, address villa, address road :  Western : Colombo : 11100"
$ws.Range("E2").Value = "123, address villa, address road :  Western : Colombo : 11100"

$ws.Range("A3").Value = "TC_SC_015 -  Choose Payment Methods"
$ws.Range("B3").Value = "test@user.com"
$ws.Range("C3").Value = 1234
$ws.Range("D3").Value = "123, address villa, address road :  Western : Colombo : 11101"
$ws.Range("E3").Value = "123, address villa, address road :  Western : Colombo : 11101"

$ws.Range("A4").Value = "TC_SC_016 - Proceed to checkout"
$ws.Range("B4").Value = "test@user.com"
$ws.Range("C4").Value = 1234
$ws.Range("D4").Value = "123, address villa, address road :  Western : Colombo : 11102"
$ws.Range("E4").Value = "123, address villa, address road :  Western : Colombo : 11102"

$ws.Range("A5").Value = "TC_SC_017 - navigating to the track orders - by clicking on the Track button of the order history page."
$ws.Range("B5").Value = "test@user.com"
$ws.Range("C5").Value = 1234
$ws.Range("D5").Value = "123, address villa, address road :  Western : Colombo : 11103"
$ws.Range("E5").Value = "123, address villa, address road :  Western : Colombo : 11103"

# "Test Status" header is added last, after the data rows (matches the
# source shared-string ordering).
$ws.Range("F1").Value = "Test Status"

# ---------------------------------------------------------------------
# 5. Hyperlinks for the "User Name" column.
# ---------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:test@user.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3:B5"), "mailto:test@user.com", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "test@user.com") | Out-Null

# ---------------------------------------------------------------------
# 6. Formatting: center + wrap text for the whole used range, row
#    heights for the data rows, and column widths.
# ---------------------------------------------------------------------
$headerRange = $ws.Range("A1:F5")
$headerRange.HorizontalAlignment = -4108
$headerRange.WrapText = $true

$ws.Rows.Item(2).RowHeight = 29
$ws.Rows.Item(3).RowHeight = 29
$ws.Rows.Item(4).RowHeight = 29
$ws.Rows.Item(5).RowHeight = 43.5

$ws.Columns.Item(1).ColumnWidth = 35.6328125
$ws.Columns.Item(2).ColumnWidth = 23.453125
$ws.Columns.Item(3).ColumnWidth = 16.26953125
$ws.Columns.Item(4).ColumnWidth = 31.7265625
$ws.Columns.Item(5).ColumnWidth = 33.7265625
$ws.Columns.Item(6).ColumnWidth = 26.36328125

# ---------------------------------------------------------------------
# 7. Final selection / active sheet state.
# ---------------------------------------------------------------------
$ws.Range("D2").Select()
$ws.Activate()
